$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recalculated timestamp precision on the existing row when the sheet was re-saved.
$ws.Range("A2").Value = 45887.66438344908

# New row 3 data: evaluation entry added by Tigist W.

$ws.Range("A3").Value = 45887.66498345505
$ws.Range("A3").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("B3").Value = "Tigist W."
$ws.Range("C3").Value = "Gemini"
$ws.Range("D3").Value = "Resistence in salmonella in ethiopia?"
$ws.Range("E3").Value = 5

$f3 = @"
---
**Human Perspective:**  
Antimicrobial resistance in Salmonella is a growing public health concern, particularly in low-resource settings like Ethiopia, where detection delays and limited treatment options can worsen outcomes (Source: aga-et-al-antibiotic-susceptibility-patterns-of-salmonella-isolates-from-clinical-food-and-environmental-sources-in (1).pdf). The situation in Ethiopia may be exacerbated by a lack of antimicrobial resistance assessments for Salmonella, weak regulations, easy access to antimicrobials without prescription, and incomplete treatment courses due to patient noncompliance (Source: IJHS-15-43.pdf). Studies in Ethiopia suggest increasing antimicrobial resistance of Salmonella to commonly used antimicrobials in the public health sector (Source: IJHS-15-43.pdf). Salmonella prevalence varies across Ethiopia: 2% in Addis Ababa (among 387 blood and stool specimens), 14% in Addis Ababa (among 387 blood specimens), 4% in Amhara (among 150 stool specimens), 3.87% in Oromia (among 232 stool specimens), 7% in Jigjiga (among 2000 stool specimens), and 1.5% in SNNP (among 381 blood specimens) (Source: IJHS-15-43.pdf). Globally, increasing multidrug resistance poses a hazard to public health (Source: IJHS-15-43.pdf).
**Animal Perspective:**  
Studies in Ethiopia suggest increasing antimicrobial resistance of Salmonella to commonly used antimicrobials in the veterinary sectors (Source: IJHS-15-43.pdf).
**Environment Perspective:**  
A study provides a comprehensive analysis of Salmonella isolates from environmental samples in Addis Ababa and nearby towns, offering current data on antimicrobial resistance patterns (Source: aga-et-al-antibiotic-susceptibility-patterns-of-salmonella-isolates-from-clinical-food-and-environmental-sources-in (1).pdf).

"@
$ws.Range("F3").Value = $f3

$g3 = @"
2.9.2 Resistance patterns in Ethiopia
Antimicrobial resistance is a global problem in general (Acha and Szyfres 2001), but it might be more severe in Ethiopia where there is lack of antimicrobial resistance assessments of Salmonella and lack of rigorous regulations but there is easy access of antimicrobials for purchase of people without prescription and incomplete treatment courses as the result of patient noncompliance (Beyene et al., 2011). There have been studies conducted in Ethiopia on salmonellosis which suggest an increase in the antimicrobial resistance of Salmonella to commonly used antimicrobials in both the public health and veterinary sectors (Mache, 2002; Molla et al., 2003; Alemayehu et al., 2004; Argaw et al., 2007; Beyene et al., 2011; Sibhat et al., 2011).
Editor Wendy A. Szymczak, Montefiore Medical Center and Albert Einstein College of Medicine, Bronx, New York, USA
IMPORTANCE Antibiotic-resistant Salmonella is a growing public health threat, particularly in low-resource settings like Ethiopia, where delayed detection and limited treatment options worsen disease outcomes. This study provides a comprehensive analysis of Salmonella isolates from clinical, food, and environmental samples in Addis Ababa and nearby towns, offering current data on antimicrobial resistance patterns. By using updated laboratory standards and sampling diverse sources, the findings highlight the urgent need for improved food safety practices, sanitation, and antimicro bial stewardship.
Address correspondence to Abebe M. Aga, agagurmu@yahoo.com, or Mesfin Tafesse Gemeda, Mesfin.tafesse@aastu.edu.et.
The authors declare no conflict of interest.
Received 31 March 2025 Accepted 23 June 2025 Published 11 July 2025
Editor Wendy A. Szymczak, Montefiore Medical Center and Albert Einstein College of Medicine, Bronx, New York, USA
IMPORTANCE Antibiotic-resistant Salmonella is a growing public health threat, particularly in low-resource settings like Ethiopia, where delayed detection and limited treatment options worsen disease outcomes. This study provides a comprehensive analysis of Salmonella isolates from clinical, food, and environmental samples in Addis Ababa and nearby towns, offering current data on antimicrobial resistance patterns. By using updated laboratory standards and sampling diverse sources, the findings highlight the urgent need for improved food safety practices, sanitation, and antimicro bial stewardship.
Address correspondence to Abebe M. Aga, agagurmu@yahoo.com, or Mesfin Tafesse Gemeda, Mesfin.tafesse@aastu.edu.et.
The authors declare no conflict of interest.
Received 31 March 2025 Accepted 23 June 2025 Published 11 July 2025
The burden of Salmonella species had different magnitudes in different parts of Ethiopia. For instance, 4% of the 150 stool specimens in Amhara [13]; 2% of among 387 blood and stool specimens in Addis Ababa [14];14% of among 387 blood specimens in Addis Ababa [15]; 3.87% of the 232 stool specimen in Oromia [16]; 7% of among 2000 stool specimen in Jigjiga [17]; and 1.5% of among 381 blood specimens in SNNP [18].
Globally, there has been an increase in multidrug resistance, which is a hazard to public
International Journal of Health Sciences
Vol. 15, Issue 1 (January - February 2021)
Abate and Assefa: Patterns of Salmonella in Ethiopia
"@
$ws.Range("G3").Value = $g3

$ws.Range("H3").Value = 5
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = 5
$ws.Range("K3").Value = ""
$ws.Range("L3").Value = 5
$ws.Range("M3").Value = ""
$ws.Range("N3").Value = 5
$ws.Range("O3").Value = ""
$ws.Range("P3").Value = ""
